$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 5 (hunk 0)
$ws.Range("H5").Value = 181.25
$ws.Range("I5").Value = 130
$ws.Range("J5").Value = 266.66666
$ws.Range("K5").Value = 130
$ws.Range("L5").Value = 266.66666
$ws.Range("M5").Value = -15
$ws.Range("N5").Value = -496.66666
# row 34 (hunk 1)
$ws.Range("H34").Value = 5100
$ws.Range("I34").Value = 5100
$ws.Range("K34").Value = 5100
$ws.Range("M34").Value = -4897
# row 36 (hunk 2)
$ws.Range("H36").Value = 5100
$ws.Range("I36").Value = 5100
$ws.Range("K36").Value = 5100
$ws.Range("M36").Value = -4385
# row 43 (hunk 3)
$ws.Range("H43").Value = 2596
$ws.Range("J43").Value = 2449.5
$ws.Range("L43").Value = 2449.5
$ws.Range("N43").Value = -2587.5
# row 62 (hunk 4)
$ws.Range("H62").Value = 4540.4
$ws.Range("I62").Value = 2925.5
$ws.Range("J62").Value = 11000
$ws.Range("K62").Value = 2925.5
$ws.Range("L62").Value = 11000
$ws.Range("M62").Value = -2301.5
$ws.Range("N62").Value = -12248
# row 65 (hunk 5)
$ws.Range("H65").Value = 4540.4
$ws.Range("I65").Value = 2925.5
$ws.Range("J65").Value = 11000
$ws.Range("K65").Value = 14627.5
$ws.Range("L65").Value = 55000
$ws.Range("M65").Value = -11507.5
$ws.Range("N65").Value = -61240
# row 70 (hunk 6)
$ws.Range("H70").Value = 7021.7144
$ws.Range("J70").Value = 10000.75
$ws.Range("L70").Value = 30002.25
$ws.Range("N70").Value = -30542.25
# row 73 (hunk 7)
$ws.Range("H73").Value = 7021.7144
$ws.Range("J73").Value = 10000.75
$ws.Range("L73").Value = 30002.25
$ws.Range("N73").Value = -31874.25
# row 97 (hunk 8)
$ws.Range("H97").Value = 732.2
$ws.Range("I97").Value = 820.6667
$ws.Range("K97").Value = 2462.0001
$ws.Range("M97").Value = -1966.0001
# row 100 (hunk 9)
$ws.Range("H100").Value = 1414.7142
$ws.Range("I100").Value = 1400.5
$ws.Range("K100").Value = 1400.5
$ws.Range("M100").Value = -859.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32 (hunk 10)
$ws.Range("H32").Value = 917
$ws.Range("I32").Value = 827.5789
$ws.Range("K32").Value = 827.5789
$ws.Range("M32").Value = -540.5789
# row 45 (hunk 11)
$ws.Range("H45").Value = 3446.5386
$ws.Range("I45").Value = 1916.5
$ws.Range("K45").Value = 1916.5
$ws.Range("M45").Value = -1539.5
# row 63 (hunk 12)
$ws.Range("H63").Value = 8999.666999999999
$ws.Range("I63").Value = 3000
$ws.Range("K63").Value = 3000
$ws.Range("M63").Value = -2314
# row 66 (hunk 13)
$ws.Range("H66").Value = 8999.666999999999
$ws.Range("I66").Value = 3000
$ws.Range("K66").Value = 15000
$ws.Range("M66").Value = -11568
# row 122 (hunk 14)
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# row 132 (hunk 15)
$ws.Range("H132").Value = 1296.4
$ws.Range("I132").Value = 1296.4
$ws.Range("K132").Value = 3889.2
$ws.Range("M132").Value = -1359.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 86 (hunk 16)
$ws.Range("H86").Value = 4008.0908
$ws.Range("I86").Value = 1618.8
$ws.Range("J86").Value = 5999.1665
$ws.Range("K86").Value = 1618.8
$ws.Range("L86").Value = 5999.1665
$ws.Range("M86").Value = -495.8
$ws.Range("N86").Value = -8245.166499999999
# row 89 (hunk 17)
$ws.Range("H89").Value = 4008.0908
$ws.Range("I89").Value = 1618.8
$ws.Range("J89").Value = 5999.1665
$ws.Range("K89").Value = 8094
$ws.Range("L89").Value = 29995.8325
$ws.Range("M89").Value = -2478
$ws.Range("N89").Value = -41227.8325
# row 99 (hunk 18)
$ws.Range("H99").Value = 1606.75
$ws.Range("I99").Value = 1606.75
$ws.Range("K99").Value = 1606.75
$ws.Range("M99").Value = -108.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 58 (hunk 19)
$ws.Range("H58").Value = 3613.9092
$ws.Range("I58").Value = 1836.5
$ws.Range("J58").Value = 5746.8
$ws.Range("K58").Value = 1836.5
$ws.Range("L58").Value = 5746.8
$ws.Range("M58").Value = -1633.5
$ws.Range("N58").Value = -6152.8
# row 62 (hunk 20)
$ws.Range("H62").Value = 2992
$ws.Range("I62").Value = 2926
$ws.Range("K62").Value = 2926
$ws.Range("M62").Value = -2302
# row 65 (hunk 21)
$ws.Range("H65").Value = 2992
$ws.Range("I65").Value = 2926
$ws.Range("K65").Value = 14630
$ws.Range("M65").Value = -11510
# row 122 (hunk 22)
$ws.Range("H122").Value = 1233.6
$ws.Range("I122").Value = 1233.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3700.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1250.8
$ws.Range("N122").ClearContents()
# row 136 (hunk 23)
$ws.Range("H136").Value = 3613.9092
$ws.Range("I136").Value = 1836.5
$ws.Range("J136").Value = 5746.8
$ws.Range("K136").Value = 5509.5
$ws.Range("L136").Value = 17240.4
$ws.Range("M136").Value = -2959.5
$ws.Range("N136").Value = -22340.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 55 (hunk 24)
$ws.Range("H55").Value = 7620
$ws.Range("J55").Value = 8542.066000000001
$ws.Range("L55").Value = 25626.198
$ws.Range("N55").Value = -25980.198
# row 80 (hunk 25)
$ws.Range("H80").Value = 4311.0835
$ws.Range("J80").Value = 5494.4
$ws.Range("L80").Value = 16483.2
$ws.Range("N80").Value = -18355.2
# row 83 (hunk 26)
$ws.Range("H83").Value = 4311.0835
$ws.Range("J83").Value = 5494.4
$ws.Range("L83").Value = 49449.6
$ws.Range("N83").Value = -58809.6
# row 86 (hunk 27)
$ws.Range("H86").Value = 558.8
$ws.Range("I86").Value = 554.2222
$ws.Range("J86").Value = 600
$ws.Range("K86").Value = 1662.6666
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -476.6666
$ws.Range("N86").Value = -4172
# row 89 (hunk 28)
$ws.Range("H89").Value = 558.8
$ws.Range("I89").Value = 554.2222
$ws.Range("J89").Value = 600
$ws.Range("K89").Value = 4987.999800000001
$ws.Range("L89").Value = 5400
$ws.Range("M89").Value = 940.0001999999995
$ws.Range("N89").Value = -17256
# row 113 (hunk 29)
$ws.Range("H113").Value = 968
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 968
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2904
$ws.Range("N113").Value = -7244
$ws.Range("M113").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 122 (hunk 30)
$ws.Range("H122").Value = 3359.1333
$ws.Range("I122").Value = 3288.5833
$ws.Range("K122").Value = 9865.749899999999
$ws.Range("M122").Value = -7415.749899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 122 (hunk 31)
$ws.Range("H122").Value = 3166.6667
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -15400
# row 136 (hunk 32)
$ws.Range("H136").Value = 1927.6666
$ws.Range("I136").Value = 1927.6666
$ws.Range("K136").Value = 5782.9998
$ws.Range("M136").Value = -3232.9998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 81 (hunk 33)
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
# row 84 (hunk 34)
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
# row 122 (hunk 35)
$ws.Range("H122").Value = 1730.909
$ws.Range("I122").Value = 1474.5
$ws.Range("K122").Value = 4423.5
$ws.Range("M122").Value = -1973.5
# row 126 (hunk 36)
$ws.Range("H126").Value = 6759.5
$ws.Range("I126").Value = 1797.5
$ws.Range("K126").Value = 5392.5
$ws.Range("M126").Value = -2922.5
